$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Report")

# Sheet1: add TANGGAL SELESAI for LOAD DATA row (shared string 36)
$ws1.Range("D2").Value = "16 Oktober 2021"

# Report sheet: add new report row - set in order D, C, B, E, A to match shared string order
$ws2.Range("D2").Value = 'Menambah tanggal pada tanggal selesai task "load data"'
$ws2.Range("C2").Value = "Selesai"
$ws2.Range("B2").Value = 'tanggal selesai pada task "load data" belum ditulis'
$ws2.Range("E2").Value = "Muhammad Fadhlan"
$ws2.Range("A2").Value = "26 Oktober 2021"

$ws2.Range("B2").WrapText = $true
$ws2.Range("D2").WrapText = $true
$ws2.Rows.Item(2).RowHeight = 43.2

# Selections
$ws1.Range("D2").Select()
$ws2.Range("C6").Select()
$ws2.Activate()
